# fall 22 week 12 complete
# Column I holds the "week 12" results; previously all marked "DNP" (Did Not
# Play) for every player, now the week is complete and the outstanding games
# are marked "NA" instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3:I10").Value = "NA"

# Leave the selection where the editor last left off while reviewing the
# updated week's column.
$ws.Range("K12").Select()
